$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Port Pin / Alternate Function" mini-table #1 (rows 19-21) ---

# Header row 19: copy header formatting (fill + border + center) from the
# existing "Sense Encode" header row (A12:B12) before writing the values.
$ws.Range("A12:B12").Copy() | Out-Null
$ws.Range("A19:B19").PasteSpecial(-4122) | Out-Null
$ws.Range("A19").Value = "Port Pin"
$ws.Range("B19").Value = "Alternate Function"

# Data row 20: A20 matches the centered/bordered style used by A13,
# B20 matches the bordered (non-center) style used by C13.
$ws.Range("A13").Copy() | Out-Null
$ws.Range("A20").PasteSpecial(-4122) | Out-Null
$ws.Range("C13").Copy() | Out-Null
$ws.Range("B20").PasteSpecial(-4122) | Out-Null
$ws.Range("A20").Value = "PA2"
$ws.Range("B20").Value = "External Interrupt"

# Spacer row 21: blank cell, centered alignment, no border/fill.
$ws.Range("A21").HorizontalAlignment = -4108
$excel.CutCopyMode = $false

# --- New "Port Pin / Alternate Function" mini-table #2 (rows 22-24) ---

$ws.Range("A12:B12").Copy() | Out-Null
$ws.Range("A22:B22").PasteSpecial(-4122) | Out-Null
$ws.Range("A22").Value = "Port Pin"
$ws.Range("B22").Value = "Alternate Function"

$ws.Range("A13").Copy() | Out-Null
$ws.Range("A23").PasteSpecial(-4122) | Out-Null
$ws.Range("C13").Copy() | Out-Null
$ws.Range("B23").PasteSpecial(-4122) | Out-Null
$ws.Range("A23").Value = "PB1"
$ws.Range("B23").Value = "Timer Limit Toggle Flag"

$ws.Range("A24").HorizontalAlignment = -4108
$excel.CutCopyMode = $false

# --- Column width tweaks ---
$ws.Columns.Item(2).ColumnWidth = 39.833333333333336
$ws.Columns.Item(5).ColumnWidth = 20.333333333333332
$ws.Columns.Item(6).ColumnWidth = 37

# --- View state: scroll down and move the selection ---
$ws.Range("C21").Select() | Out-Null
